# Refresh the movie report table (Title/Genre/Rating/Year) with the new
# dataset. Only cells whose value actually changes vs. the original are
# written. Rating and Year are stored as text in this report, so numeric-
# looking values are written with a leading apostrophe to keep Excel from
# auto-converting them to the Number type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Mary"
$ws.Cells.Item(2, 2).Value = "Action"
$ws.Cells.Item(2, 3).Value = "'5.2"

$ws.Cells.Item(3, 1).Value = "Red One"
$ws.Cells.Item(3, 2).Value = "Action"
$ws.Cells.Item(3, 3).Value = "'6.7"

$ws.Cells.Item(4, 1).Value = "Heretic"
$ws.Cells.Item(4, 3).Value = "'7.1"

$ws.Cells.Item(5, 1).Value = "Venom: The Last Dance"
$ws.Cells.Item(5, 2).Value = "Action"
$ws.Cells.Item(5, 3).Value = "'6.1"

$ws.Cells.Item(6, 1).Value = "[TA] Thangalaan"
$ws.Cells.Item(6, 2).Value = "Action"
$ws.Cells.Item(6, 3).Value = "'6.9"

$ws.Cells.Item(7, 1).Value = "The Best Christmas Pageant Ever"
$ws.Cells.Item(7, 2).Value = "Adventure"
$ws.Cells.Item(7, 3).Value = "'7"

$ws.Cells.Item(8, 1).Value = "Juror #2"
$ws.Cells.Item(8, 2).Value = "Crime"
$ws.Cells.Item(8, 3).Value = "'7.1"

$ws.Cells.Item(9, 1).Value = "That Christmas"
$ws.Cells.Item(9, 2).Value = "Adventure"
$ws.Cells.Item(9, 3).Value = "'6.8"

$ws.Cells.Item(10, 1).Value = "Holiday Touchdown: A Chiefs Love Story"
$ws.Cells.Item(10, 2).Value = "Comedy"
$ws.Cells.Item(10, 3).Value = "'6.4"

$ws.Cells.Item(11, 1).Value = "[HI] Sikandar Ka Muqaddar"
$ws.Cells.Item(11, 2).Value = "Action"
$ws.Cells.Item(11, 3).Value = "'6.1"

$ws.Cells.Item(12, 1).Value = "[TE] Lucky Baskhar"
$ws.Cells.Item(12, 2).Value = "Crime"
$ws.Cells.Item(12, 3).Value = "'8.1"
$ws.Cells.Item(12, 4).Value = "'2024"

$ws.Cells.Item(13, 1).Value = "[FR] The Seed of the Sacred Fig"
$ws.Cells.Item(13, 2).Value = "Crime"
$ws.Cells.Item(13, 3).Value = "'7.6"
$ws.Cells.Item(13, 4).Value = "'2024"

$ws.Cells.Item(14, 1).Value = "The Convert"
$ws.Cells.Item(14, 3).Value = "'6.4"
$ws.Cells.Item(14, 4).Value = "'2023"

$ws.Cells.Item(15, 1).Value = "[KO] Project Silence"
$ws.Cells.Item(15, 2).Value = "Action"
$ws.Cells.Item(15, 3).Value = "'5.5"
$ws.Cells.Item(15, 4).Value = "'2023"

$ws.Cells.Item(16, 1).Value = "David Attenborough: A Life on Our Planet"
$ws.Cells.Item(16, 2).Value = "Action"
$ws.Cells.Item(16, 3).Value = "'8.9"
$ws.Cells.Item(16, 4).Value = "'2020"

$ws.Cells.Item(17, 1).Value = "The Substance"
$ws.Cells.Item(17, 2).Value = "Drama"
$ws.Cells.Item(17, 3).Value = "'7.4"
$ws.Cells.Item(17, 4).Value = "'2024"

$ws.Cells.Item(18, 1).Value = "[HI] Jigra"
$ws.Cells.Item(18, 2).Value = "Action"
$ws.Cells.Item(18, 3).Value = "'6.4"

$ws.Cells.Item(19, 1).Value = "[FR] Cat and Dog"
$ws.Cells.Item(19, 3).Value = "'5.2"
$ws.Cells.Item(19, 4).Value = "'2024"

$ws.Cells.Item(20, 1).Value = "Conclave"
$ws.Cells.Item(20, 2).Value = "Drama"
$ws.Cells.Item(20, 3).Value = "'7.4"
$ws.Cells.Item(20, 4).Value = "'2024"

$ws.Cells.Item(21, 1).Value = "Freud's Last Session"
$ws.Cells.Item(21, 3).Value = "'6.1"
